$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 78032.766
$ws.Range("I61").Value = 712.1111
$ws.Range("J61").Value = 252004.25
$ws.Range("K61").Value = 2136.3333
$ws.Range("L61").Value = 756012.75
$ws.Range("M61").Value = -1964.3333
$ws.Range("N61").Value = -756356.75
$ws.Range("H69").Value = 3846
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 4410
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 13230
$ws.Range("M69").Value = -8126
$ws.Range("N69").Value = -14978
$ws.Range("H72").Value = 3846
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 4410
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 39690
$ws.Range("M72").Value = -22632
$ws.Range("N72").Value = -48426
$ws.Range("H98").Value = 9334.223
$ws.Range("I98").Value = 9001.5
$ws.Range("J98").Value = 9999.666999999999
$ws.Range("K98").Value = 9001.5
$ws.Range("L98").Value = 9999.666999999999
$ws.Range("M98").Value = -7503.5
$ws.Range("N98").Value = -12995.667
$ws.Range("H112").Value = 19308732
$ws.Range("I112").Value = 850
$ws.Range("J112").Value = 27890012
$ws.Range("K112").Value = 2550
$ws.Range("L112").Value = 83670036
$ws.Range("M112").Value = -1442
$ws.Range("N112").Value = -83672252
$ws.Range("H118").Value = 917.58826
$ws.Range("I118").Value = 364.875
$ws.Range("J118").Value = 1408.8889
$ws.Range("K118").Value = 1094.625
$ws.Range("L118").Value = 4226.6667
$ws.Range("M118").Value = 562.375
$ws.Range("N118").Value = -7540.6667
$ws.Range("H122").Value = 9334.223
$ws.Range("I122").Value = 9001.5
$ws.Range("J122").Value = 9999.666999999999
$ws.Range("K122").Value = 27004.5
$ws.Range("L122").Value = 29999.001
$ws.Range("M122").Value = -24554.5
$ws.Range("N122").Value = -34899.001
$ws.Range("H129").Value = 4167776
$ws.Range("I129").Value = 83335464
$ws.Range("K129").Value = 250006392
$ws.Range("M129").Value = -250001392

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 45864.145
$ws.Range("I21").Value = 5249.5
$ws.Range("K21").Value = 5249.5
$ws.Range("M21").Value = -4875.5
$ws.Range("H32").Value = 3676.6
$ws.Range("I32").Value = 3676.6
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3676.6
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3389.6
$ws.Range("N32").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H97").Value = 529.4
$ws.Range("I97").Value = 476
$ws.Range("J97").Value = 1010
$ws.Range("K97").Value = 476
$ws.Range("L97").Value = 1010
$ws.Range("M97").Value = 20
$ws.Range("N97").Value = -2002
$ws.Range("H122").Value = 2565.889
$ws.Range("I122").Value = 1969.9584
$ws.Range("K122").Value = 5909.8752
$ws.Range("M122").Value = -3459.8752

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3000.8125
$ws.Range("I20").Value = 2789.125
$ws.Range("J20").Value = 3212.5
$ws.Range("K20").Value = 2789.125
$ws.Range("L20").Value = 3212.5
$ws.Range("M20").Value = -2542.125
$ws.Range("N20").Value = -3706.5
$ws.Range("H23").Value = 58691.715
$ws.Range("I23").Value = 16933.334
$ws.Range("J23").Value = 90010.5
$ws.Range("K23").Value = 16933.334
$ws.Range("L23").Value = 90010.5
$ws.Range("M23").Value = -16650.334
$ws.Range("N23").Value = -90576.5
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622
$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112
$ws.Range("H132").Value = 29833.334
$ws.Range("J132").Value = 29833.334
$ws.Range("L132").Value = 29833.334
$ws.Range("N132").Value = -39953.334

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 27085.7
$ws.Range("I39").Value = 9000
$ws.Range("J39").Value = 39142.832
$ws.Range("K39").Value = 9000
$ws.Range("L39").Value = 39142.832
$ws.Range("M39").Value = -8609
$ws.Range("N39").Value = -39924.832
$ws.Range("H49").Value = 27085.7
$ws.Range("I49").Value = 9000
$ws.Range("J49").Value = 39142.832
$ws.Range("K49").Value = 9000
$ws.Range("L49").Value = 39142.832
$ws.Range("M49").Value = -8818
$ws.Range("N49").Value = -39506.832
$ws.Range("H68").Value = 29929.834
$ws.Range("J68").Value = 29929.834
$ws.Range("L68").Value = 29929.834
$ws.Range("N68").Value = -31427.834
$ws.Range("H69").Value = 40201
$ws.Range("J69").Value = 40201
$ws.Range("L69").Value = 40201
$ws.Range("N69").Value = -41699
$ws.Range("H71").Value = 29929.834
$ws.Range("J71").Value = 29929.834
$ws.Range("L71").Value = 89789.50199999999
$ws.Range("N71").Value = -97277.50199999999
$ws.Range("H72").Value = 40201
$ws.Range("J72").Value = 40201
$ws.Range("L72").Value = 120603
$ws.Range("N72").Value = -128091
$ws.Range("H105").Value = 3451.2
$ws.Range("I105").Value = 5502.25
$ws.Range("J105").Value = 2705.3635
$ws.Range("K105").Value = 5502.25
$ws.Range("L105").Value = 2705.3635
$ws.Range("M105").Value = -3755.25
$ws.Range("N105").Value = -6199.363499999999
$ws.Range("H122").Value = 2368.303
$ws.Range("I122").Value = 1971.1923
$ws.Range("J122").Value = 3843.2856
$ws.Range("K122").Value = 5913.5769
$ws.Range("L122").Value = 11529.8568
$ws.Range("M122").Value = -3463.5769
$ws.Range("N122").Value = -16429.8568

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 590.2
$ws.Range("I114").Value = 188.23077
$ws.Range("J114").Value = 897.58826
$ws.Range("K114").Value = 564.69231
$ws.Range("L114").Value = 2692.76478
$ws.Range("M114").Value = 2689.30769
$ws.Range("N114").Value = -9200.76478
$ws.Range("H116").Value = 1239.7778
$ws.Range("I116").Value = 159.66667
$ws.Range("K116").Value = 479.00001
$ws.Range("M116").Value = 2962.99999
$ws.Range("H131").Value = 1135.4062
$ws.Range("J131").Value = 1077.0182
$ws.Range("L131").Value = 3231.0546
$ws.Range("N131").Value = -13311.0546

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4816
$ws.Range("I70").Value = 4778.3335
$ws.Range("J70").Value = 4966.6665
$ws.Range("K70").Value = 4778.3335
$ws.Range("L70").Value = 4966.6665
$ws.Range("M70").Value = -4508.3335
$ws.Range("N70").Value = -5506.6665
$ws.Range("H73").Value = 4816
$ws.Range("I73").Value = 4778.3335
$ws.Range("J73").Value = 4966.6665
$ws.Range("K73").Value = 4778.3335
$ws.Range("L73").Value = 4966.6665
$ws.Range("M73").Value = -3842.3335
$ws.Range("N73").Value = -6838.6665

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 34000
$ws.Range("J36").Value = 34000
$ws.Range("L36").Value = 34000
$ws.Range("N36").Value = -35124
$ws.Range("H59").Value = 20199
$ws.Range("J59").Value = 20199
$ws.Range("L59").Value = 20199
$ws.Range("N59").Value = -21507
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622
$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112
$ws.Range("H86").Value = 29700
$ws.Range("J86").Value = 29700
$ws.Range("L86").Value = 29700
$ws.Range("N86").Value = -32072
$ws.Range("H89").Value = 29700
$ws.Range("J89").Value = 29700
$ws.Range("L89").Value = 89100
$ws.Range("N89").Value = -100956
$ws.Range("H132").Value = 2134.6726
$ws.Range("I132").Value = 1425.8649
$ws.Range("K132").Value = 4277.5947
$ws.Range("M132").Value = -1747.5947

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 13648.429
$ws.Range("J54").Value = 13648.429
$ws.Range("L54").Value = 13648.429
$ws.Range("N54").Value = -14688.429
$ws.Range("H62").Value = 3950.3333
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 4040.4
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 4040.4
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -5288.4
$ws.Range("H65").Value = 3950.3333
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 4040.4
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 20202
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -26442
$ws.Range("H75").Value = 37457.5
$ws.Range("J75").Value = 37457.5
$ws.Range("L75").Value = 37457.5
$ws.Range("N75").Value = -39329.5
$ws.Range("H78").Value = 37457.5
$ws.Range("J78").Value = 37457.5
$ws.Range("L78").Value = 112372.5
$ws.Range("N78").Value = -121732.5
$ws.Range("H122").Value = 2074.8386
$ws.Range("I122").Value = 1677.7307
$ws.Range("J122").Value = 4139.8
$ws.Range("K122").Value = 5033.1921
$ws.Range("L122").Value = 12419.4
$ws.Range("M122").Value = -2583.1921
$ws.Range("N122").Value = -17319.4
